$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update name, email, and repo link in row 2
$ws.Range("A2").Value = "محمد زكي جلال"
$ws.Range("B2").Value = "mohammedzakigalal@gmail.com"
$ws.Range("C2").Value = "https://github.com/Iammohamedzaki/Security-Task"

# Update the active selection to A3
$ws.Range("A3").Select()
